# Mises a jour du 11 juillet
# Responsables d'EC et coefficients en ME

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level: rename sheet "ME" -> "MANE"
# ---------------------------------------------------------------------------
$wsGEEL = $wb.Worksheets.Item(1)
$wsGGL  = $wb.Worksheets.Item(2)
$wsGAGL = $wb.Worksheets.Item(3)
$wsMANE = $wb.Worksheets.Item(4)
$wsMANE.Name = "MANE"

# ---------------------------------------------------------------------------
# 2) GEEL sheet: view state only (selection moved, no longer the active tab)
# ---------------------------------------------------------------------------
$wsGEEL.Activate()
$excel.ActiveWindow.ScrollRow = 58
$wsGEEL.Range("D60").Select()

# ---------------------------------------------------------------------------
# 3) GGL sheet: view state only (scrolled)
# ---------------------------------------------------------------------------
$wsGGL.Activate()
$excel.ActiveWindow.ScrollRow = 30
$wsGGL.Range("N38").Select()

# ---------------------------------------------------------------------------
# 4) GAGL sheet: view state only (scrolled + new selection)
# ---------------------------------------------------------------------------
$wsGAGL.Activate()
$excel.ActiveWindow.ScrollRow = 38
$wsGAGL.Range("H59").Select()

# ---------------------------------------------------------------------------
# 5) MANE sheet (formerly "ME"): content + formatting updates, becomes the
#    active / selected tab.
# ---------------------------------------------------------------------------

# --- Responsables d'EC (column F) ------------------------------------------
# GSIC-20204A: responsable supprime (cellule videe)
$wsMANE.Range("F31").Value = " "

# STAG-O-20228A / MIR-20211A: JF. Berthevas -> F. Mayon
$wsMANE.Range("F34").Value = "F. Mayon"
$wsMANE.Range("F35").Value = "F. Mayon"

# PROJ-20312C (gestion de projet) + Stage/Memoire GEEL-recherche: J. Viau / F. Mayon -> E. Lamendour
$wsMANE.Range("F56").Value = "E. Lamendour"
$wsMANE.Range("F62").Value = "E. Lamendour"
$wsMANE.Range("F63").Value = "E. Lamendour"

# --- Coefficients (column H) ------------------------------------------------
# UE rows whose coefficient is cleared entirely (no value)
$wsMANE.Range("H7").ClearContents()
$wsMANE.Range("H11").ClearContents()
$wsMANE.Range("H15").ClearContents()
$wsMANE.Range("H22").ClearContents()
$wsMANE.Range("H26").ClearContents()

# EC rows (under UE-20201C) that now get a coefficient of 1, with new
# formatting: 8pt red "Athelas Regular", centered horizontally.
foreach ($addr in @("H23", "H24", "H25")) {
    $c = $wsMANE.Range($addr)
    $c.Value = 1
    $c.Font.Name = "Athelas Regular"
    $c.Font.Size = 8
    $c.Font.Color = 255
    $c.HorizontalAlignment = -4108
}

# UE-20211A / STAG-O-20228A / MIR-20211A rows: coefficient 3 -> 1, with new
# formatting: 8pt red "Athelas Regular", centered horizontally + vertically.
# H33 additionally gets the grey (D9D9D9) fill (UE-level row), H34/H35 (EC
# rows) keep no fill.
$c = $wsMANE.Range("H33")
$c.Value = 1
$c.Font.Name = "Athelas Regular"
$c.Font.Size = 8
$c.Font.Color = 255
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$c.Interior.Color = 14277081

foreach ($addr in @("H34", "H35")) {
    $c = $wsMANE.Range($addr)
    $c.Value = 1
    $c.Font.Name = "Athelas Regular"
    $c.Font.Size = 8
    $c.Font.Color = 255
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# UE coefficient cells whose numeric coefficient is replaced by a blank
# (single-space) placeholder text instead of a number.
foreach ($addr in @("H30", "H38", "H41", "H43", "H50")) {
    $wsMANE.Range($addr).Value = " "
}

# Same blank placeholder, but also restyled (grey fill + centered +
# vcentered, 8pt red Athelas Regular) - UE-20303C / UE-20305C rows.
foreach ($addr in @("H47", "H54")) {
    $c = $wsMANE.Range($addr)
    $c.Value = " "
    $c.Font.Name = "Athelas Regular"
    $c.Font.Size = 8
    $c.Font.Color = 255
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.Interior.Color = 14277081
}

# UE-20306C / MS-20311C (H48 / H55): coefficient 2 -> 1 with the same
# centered/vcentered red-Athelas formatting (no fill).
foreach ($addr in @("H48", "H55")) {
    $c = $wsMANE.Range($addr)
    $c.Value = 1
    $c.Font.Name = "Athelas Regular"
    $c.Font.Size = 8
    $c.Font.Color = 255
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# --- View state: MANE becomes the selected/active tab ----------------------
$wsMANE.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$wsMANE.Range("H33:H35").Select()
